$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct status name/labels in the "statut_label" (B) and "statut_name" (C) columns.

$used = $ws.UsedRange
$rows = $used.Rows.Count

for ($r = 2; $r -le $rows; $r++) {
    $label = $ws.Cells.Item($r, 2).Value2
    $name  = $ws.Cells.Item($r, 3).Value2

    if ($label -eq "bleu") {
        $ws.Cells.Item($r, 2).Value2 = "noir"
    }

    if ($name -eq "pas de résultat ni de publication") {
        $ws.Cells.Item($r, 3).Value2 = "pas de résultat postés ni publiés"
    }
    elseif ($name -eq "résultat et / ou publication posté dans les 36 mois") {
        $ws.Cells.Item($r, 3).Value2 = "résultat postés ou publiés dans les 36 mois"
    }
    elseif ($name -eq "résultat et / ou publication posté") {
        $ws.Cells.Item($r, 3).Value2 = "résultat postés ou publiés"
    }
}
